# Update "想去人数" (F column) counts on 展览, 演出, and 全部类型 sheets
# to reflect output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4714
$ws1.Range("F3").Value = 1876
$ws1.Range("F6").Value = 3175
$ws1.Range("F9").Value = 284
$ws1.Range("F10").Value = 649
$ws1.Range("F12").Value = 546
$ws1.Range("F14").Value = 141
$ws1.Range("F15").Value = 1795
$ws1.Range("F16").Value = 1384
$ws1.Range("F18").Value = 1648
$ws1.Range("F22").Value = 19
$ws1.Range("F24").Value = 542
$ws1.Range("F27").Value = 112
$ws1.Range("F30").Value = 45
$ws1.Range("F32").Value = 4004
$ws1.Range("F33").Value = 11
$ws1.Range("F34").Value = 783
$ws1.Range("F36").Value = 1543
$ws1.Range("F38").Value = 1899

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 58

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4714
$ws4.Range("F3").Value = 1876
$ws4.Range("F6").Value = 3175
$ws4.Range("F9").Value = 284
$ws4.Range("F10").Value = 649
$ws4.Range("F12").Value = 546
$ws4.Range("F15").Value = 141
$ws4.Range("F16").Value = 1795
$ws4.Range("F17").Value = 1384
$ws4.Range("F19").Value = 1648
$ws4.Range("F23").Value = 19
$ws4.Range("F25").Value = 542
$ws4.Range("F28").Value = 112
$ws4.Range("F31").Value = 45
$ws4.Range("F33").Value = 4004
$ws4.Range("F34").Value = 58
$ws4.Range("F35").Value = 11
$ws4.Range("F37").Value = 783
$ws4.Range("F39").Value = 1543
$ws4.Range("F41").Value = 1899
